$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "test3@gmail.com"
$ws.Range("B6").Value = "blacklisted"
